$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E3").Value = 2
